$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to c160a3b9-9890-481e-96a8-76025723d394.md
# Mark it ready for handoff in both locale columns.
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the same file; update status + new handoff datetime.
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-01-11 16:38:08"

# de-de sheet: row 3 is the same file; update status + new handoff datetime.
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-01-11 16:38:35"
